$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1001
$ws.Range("F6").Value = 7244
$ws.Range("F8").Value = 922
$ws.Range("F9").Value = 274
$ws.Range("F10").Value = 766
$ws.Range("F11").Value = 550
$ws.Range("F12").Value = 65
$ws.Range("F13").Value = 58
$ws.Range("F15").Value = 837
$ws.Range("F16").Value = 2884
$ws.Range("F17").Value = 153
$ws.Range("F18").Value = 41
$ws.Range("F19").Value = 261
$ws.Range("F20").Value = 737
$ws.Range("F22").Value = 428
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 145
$ws.Range("F25").Value = 188
$ws.Range("F26").Value = 134
$ws.Range("F27").Value = 195
$ws.Range("F29").Value = 70
$ws.Range("F30").Value = 183
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 301
$ws.Range("F34").Value = 383
$ws.Range("F38").Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 195

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 195
$ws.Range("F3").Value = 1001
$ws.Range("F10").Value = 7244
$ws.Range("F12").Value = 922
$ws.Range("F13").Value = 274
$ws.Range("F14").Value = 766
$ws.Range("F15").Value = 550
$ws.Range("F16").Value = 65
$ws.Range("F17").Value = 58
$ws.Range("F19").Value = 837
$ws.Range("F21").Value = 2884
$ws.Range("F22").Value = 153
$ws.Range("F23").Value = 41
$ws.Range("F25").Value = 261
$ws.Range("F26").Value = 737
$ws.Range("F29").Value = 428
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 145
$ws.Range("F32").Value = 188
$ws.Range("F33").Value = 134
$ws.Range("F34").Value = 195
$ws.Range("F36").Value = 70
$ws.Range("F37").Value = 183
$ws.Range("F39").Value = 17
$ws.Range("F40").Value = 301
$ws.Range("F41").Value = 383
$ws.Range("F45").Value = 40
